$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the columns we touch so numeric-looking strings
# (e.g. "0.9978", "1.001") are preserved as text, matching the source
# workbook where these are inlineStr cells, not numbers.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '24.921.42'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -3.84%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.636.03'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -6.04%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9978'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.35%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '235.56'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -4.37%  '
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.02%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4730'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -6.51%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2555'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -5.64%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06001'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -2.81%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07165'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -0.91%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.638.08'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -5.96%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '14.73'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -2.22%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.6146'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -4.85%  '
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -4.62%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '72.48'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -6.42%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.001'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +0.01%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.9975'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -0.39%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '24.913.78'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000006561'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -3.43%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.20'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -5.08%  '
$ws.Range('B21').NumberFormat = '@'
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').NumberFormat = '@'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.403'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +2.97%  '
$ws.Range('B22').NumberFormat = '@'
$ws.Range('B22').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C22').NumberFormat = '@'
$ws.Range('C22').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.846.42'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -6.09%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.563'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -0.69%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.255'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -2.18%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '132.44'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -2.69%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '14.77'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -2.77%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.369'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -8.77%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '102.39'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -2.74%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.652'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -6.54%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.718'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -4.68%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.07746'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -5.71%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -2.66%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.04373'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -6.34%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.9991'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -0.10%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -2.09%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.9155'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -7.60%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.5799'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -6.51%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.533'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -6.97%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01553'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -2.67%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.9983'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -0.22%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8235'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +9.04%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.792'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -6.10%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '97.32'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -1.73%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.3705'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -3.58%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '4.740'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -4.70%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.1133'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +0.13%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '6.074'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -2.80%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.05193'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -0.79%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '29.50'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -3.68%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.000'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -0.24%  '
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -0.46%  '
